$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.1990728830790331
$ws.Range("D2").Value = 0.8440340933708055

$ws.Range("C3").Value = 0.5199840443020325
$ws.Range("D3").Value = 0.6082679190106148

$ws.Range("C4").Value = -0.3250140634376487
$ws.Range("D4").Value = 0.7482409815123474

$ws.Range("C5").Value = 1.530271258711054
$ws.Range("D5").Value = 0.1402027967931339

$ws.Range("C6").Value = 0.7664205132818023
$ws.Range("D6").Value = 0.4515722526766579

$ws.Range("C7").Value = -0.1455810213335589
$ws.Range("D7").Value = 0.8855777741370552

$ws.Range("C8").Value = 1.696411487473054
$ws.Range("D8").Value = 0.1039140238237213

$ws.Range("C9").Value = -1.045130898609187
$ws.Range("D9").Value = 0.3073107159754265

$ws.Range("C10").Value = 1.03341455275409
$ws.Range("D10").Value = 0.3126370607486242

$ws.Range("C11").Value = 1.658375114688425
$ws.Range("D11").Value = 0.1114331761163869
